$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared text "E7420" -> "E7420L" (affects G2:G27 which all share this string)
$ws.Range("G2:G27").Value = "E7420L"

# Replace the =FALSE() formula cells in H2:H27 with a literal boolean FALSE value
$ws.Range("H2:H27").Value = $false
